$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.776.13'
$ws.Range('E2').Value = '  -2.67%  '
$ws.Range('D3').Value = '1.786.74'
$ws.Range('E3').Value = '  -2.24%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.88'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.10%  '
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5109'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.52%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3865'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.53%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07828'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -6.83%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.090'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '40.72'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.83%  '
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.219'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -3.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.18'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -4.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.220'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -4.04%  '
$ws.Range('D16').Value = '1.773.99'
$ws.Range('E16').Value = '  -2.62%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '91.23'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001075'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -4.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06529'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.21%  '
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.01'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -4.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.904'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.75%  '
$ws.Range('D23').Value = '27.828.39'
$ws.Range('E23').Value = '  -2.55%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.01'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.222'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.63'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.69%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.21'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -4.30%  '
$ws.Range('D28').Value = '1.987.67'
$ws.Range('E28').Value = '  -2.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.358'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '123.54'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.78%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1074'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.81%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.034'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -6.04%  '
$ws.Range('E33').Value = '  -0.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.479'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -4.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07080'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -5.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.740'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2122'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -4.44%  '
$ws.Range('E39').Value = '  -0.35%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.978'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -4.79%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6075'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -4.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.145'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -3.91%  '
$ws.Range('B44').Value = 'WEMIXTOKEN'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.313'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -6.27%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.10'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -3.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5881'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.692'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.07'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.33%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.195'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.70%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.900'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -4.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06812'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.47%  '
